$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027637184202276
$ws.Range("D2").Value = 1.035249955163783
$ws.Range("E2").Value = 1.027688866148185
$ws.Range("F2").Value = 1.043871701773545
$ws.Range("I2").Value = 1.031462112632618
$ws.Range("J2").Value = 1.032793951103716
$ws.Range("K2").Value = 1.038047099433764
$ws.Range("L2").Value = 1.03050786416438
$ws.Range("M2").Value = 1.046644346015466
$ws.Range("N2").Value = 1.014879934026104
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028784106175942
$ws.Range("D3").Value = 1.036119350626239
$ws.Range("E3").Value = 1.028668675328348
$ws.Range("F3").Value = 1.044977917598278
$ws.Range("I3").Value = 1.031653437355852
$ws.Range("J3").Value = 1.033580128986484
$ws.Range("K3").Value = 1.038725691147029
$ws.Range("L3").Value = 1.031294986460946
$ws.Range("M3").Value = 1.047560914134618
$ws.Range("N3").Value = 1.015146306774611
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029526246915702
$ws.Range("D4").Value = 1.036681712935732
$ws.Range("E4").Value = 1.029303023986771
$ws.Range("F4").Value = 1.045693995596265
$ws.Range("I4").Value = 1.031775765471209
$ws.Range("J4").Value = 1.034088336042094
$ws.Range("K4").Value = 1.039163954873076
$ws.Range("L4").Value = 1.031804043465944
$ws.Range("M4").Value = 1.048153689856325
$ws.Range("N4").Value = 1.015318341724267
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.029838245247344
$ws.Range("D5").Value = 1.036918083358537
$ws.Range("E5").Value = 1.029569787296925
$ws.Range("F5").Value = 1.045995102660339
$ws.Range("I5").Value = 1.031826839771107
$ws.Range("J5").Value = 1.034301866262728
$ws.Range("K5").Value = 1.039348002186363
$ws.Range("L5").Value = 1.032017988099047
$ws.Range("M5").Value = 1.048402819909064
$ws.Range("N5").Value = 1.015390587167158
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.029890631309066
$ws.Range("D6").Value = 1.036957768252117
$ws.Range("E6").Value = 1.029614582910114
$ws.Range("F6").Value = 1.046045663849327
$ws.Range("I6").Value = 1.031835394700822
$ws.Range("J6").Value = 1.034337711896693
$ws.Range("K6").Value = 1.039378892886729
$ws.Range("L6").Value = 1.032053906660144
$ws.Range("M6").Value = 1.048444645674106
$ws.Range("N6").Value = 1.015402712918087
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02953041584208
$ws.Range("D7").Value = 1.036684871512471
$ws.Range("E7").Value = 1.029306588164323
$ws.Range("F7").Value = 1.04569801873656
$ws.Range("I7").Value = 1.031776449313471
$ws.Range("J7").Value = 1.034091189714602
$ws.Range("K7").Value = 1.039166414903193
$ws.Range("L7").Value = 1.031806902451622
$ws.Range("M7").Value = 1.048157019030343
$ws.Range("N7").Value = 1.015319307378092
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028024791333814
$ws.Range("D8").Value = 1.035543811693491
$ws.Range("E8").Value = 1.02801992557716
$ws.Range("F8").Value = 1.044245494282446
$ws.Range("I8").Value = 1.031527076402764
$ws.Range("J8").Value = 1.033059747847973
$ws.Range("K8").Value = 1.038276604677311
$ws.Range("L8").Value = 1.03077393056953
$ws.Range("M8").Value = 1.046954167798936
$ws.Range("N8").Value = 1.014970023452184
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02537168405284
$ws.Range("D9").Value = 1.033531630429703
$ws.Range("E9").Value = 1.025755308641166
$ws.Range("F9").Value = 1.041688110721268
$ws.Range("I9").Value = 1.031076381353122
$ws.Range("J9").Value = 1.031238352280681
$ws.Range("K9").Value = 1.036702281427355
$ws.Range("L9").Value = 1.028951672132878
$ws.Range("M9").Value = 1.044832244119104
$ws.Range("N9").Value = 1.014352043194359
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023602874717165
$ws.Range("D10").Value = 1.032189176832775
$ws.Range("E10").Value = 1.024247330083994
$ws.Range("F10").Value = 1.039984613065469
$ws.Range("I10").Value = 1.030768346054885
$ws.Range("J10").Value = 1.030021469276079
$ws.Range("K10").Value = 1.035648446512734
$ws.Range("L10").Value = 1.027735451885539
$ws.Range("M10").Value = 1.043416035235867
$ws.Range("N10").Value = 1.013938375723403
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022836927908751
$ws.Range("D11").Value = 1.031607641809122
$ws.Range("E11").Value = 1.023594773666733
$ws.Range("F11").Value = 1.039247312040967
$ws.Range("I11").Value = 1.030633167882632
$ws.Range("J11").Value = 1.02949391857699
$ws.Range("K11").Value = 1.035191105311086
$ws.Range("L11").Value = 1.027208483021264
$ws.Range("M11").Value = 1.04280241886905
$ws.Range("N11").Value = 1.013758853836217
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022552413569436
$ws.Range("D12").Value = 1.031391597010488
$ws.Range("E12").Value = 1.02335244605182
$ws.Range("F12").Value = 1.038973493793444
$ws.Range("I12").Value = 1.030582686709074
$ws.Range("J12").Value = 1.029297867172678
$ws.Range("K12").Value = 1.035021074309534
$ws.Range("L12").Value = 1.027012692143489
$ws.Range("M12").Value = 1.042574435558429
$ws.Range("N12").Value = 1.013692111047251
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022613443228252
$ws.Range("D13").Value = 1.03143794102123
$ws.Range("E13").Value = 1.023404423399371
$ws.Range("F13").Value = 1.03903222656641
$ws.Range("I13").Value = 1.030593527309797
$ws.Range("J13").Value = 1.029339925199277
$ws.Range("K13").Value = 1.03505755351567
$ws.Range("L13").Value = 1.027054692275006
$ws.Range("M13").Value = 1.042623341424406
$ws.Range("N13").Value = 1.013706430335391
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022813410028103
$ws.Range("D14").Value = 1.031589784230371
$ws.Range("E14").Value = 1.02357474155084
$ws.Range("F14").Value = 1.03922467715981
$ws.Range("I14").Value = 1.030629000600762
$ws.Range("J14").Value = 1.029477714860057
$ws.Range("K14").Value = 1.035177053641193
$ws.Range("L14").Value = 1.027192299933017
$ws.Range("M14").Value = 1.04278357490127
$ws.Range("N14").Value = 1.013753338088508
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022936615123179
$ws.Range("D15").Value = 1.031683334961561
$ws.Range("E15").Value = 1.023679688240914
$ws.Range("F15").Value = 1.039343258688544
$ws.Range("I15").Value = 1.030650821090179
$ws.Range("J15").Value = 1.029562598934156
$ws.Range("K15").Value = 1.035250661181265
$ws.Range("L15").Value = 1.027277077756823
$ws.Range("M15").Value = 1.042882292219056
$ws.Range("N15").Value = 1.013782231501156
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023653707100804
$ws.Range("D16").Value = 1.032227766246403
$ws.Range("E16").Value = 1.02429064666799
$ws.Range("F16").Value = 1.040033552096017
$ws.Range("I16").Value = 1.030777279525229
$ws.Range("J16").Value = 1.030056467706181
$ws.Range("K16").Value = 1.035678777146652
$ws.Range("L16").Value = 1.02777041792272
$ws.Range("M16").Value = 1.043456750704478
$ws.Range("N16").Value = 1.013950281539719
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024103508015962
$ws.Range("D17").Value = 1.032569208223873
$ws.Range("E17").Value = 1.0246739937374
$ws.Range("F17").Value = 1.040466641421445
$ws.Range("I17").Value = 1.030856122474378
$ws.Range("J17").Value = 1.030366088834505
$ws.Range("K17").Value = 1.035947048680067
$ws.Range("L17").Value = 1.028079786488364
$ws.Range("M17").Value = 1.043816988956656
$ws.Range("N17").Value = 1.014055587421943
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024365865452501
$ws.Range("D18").Value = 1.032768342045564
$ws.Range("E18").Value = 1.024897633017144
$ws.Range("F18").Value = 1.040719286467161
$ws.Range("I18").Value = 1.030901936863869
$ws.Range("J18").Value = 1.030546624629175
$ws.Range("K18").Value = 1.036103428172749
$ws.Range("L18").Value = 1.028260203409625
$ws.Range("M18").Value = 1.044027072375667
$ws.Range("N18").Value = 1.014116971863746
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024455321954194
$ws.Range("D19").Value = 1.032836237550787
$ws.Range("E19").Value = 1.024973894931855
$ws.Range("F19").Value = 1.040805437237889
$ws.Range("I19").Value = 1.030917528981108
$ws.Range("J19").Value = 1.030608172345715
$ws.Range("K19").Value = 1.03615673279763
$ws.Range("L19").Value = 1.028321715414385
$ws.Range("M19").Value = 1.044098699053784
$ws.Range("N19").Value = 1.014137895811303
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024055249040953
$ws.Range("D20").Value = 1.032532577163642
$ws.Range("E20").Value = 1.024632860150467
$ws.Range("F20").Value = 1.04042017179502
$ws.Range("I20").Value = 1.030847681304571
$ws.Range("J20").Value = 1.030332875728269
$ws.Range("K20").Value = 1.035918275900141
$ws.Range("L20").Value = 1.028046597537037
$ws.Range("M20").Value = 1.043778342679386
$ws.Range("N20").Value = 1.014044293100737
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022754525010997
$ws.Range("D21").Value = 1.031545071207108
$ws.Range("E21").Value = 1.023524585422897
$ws.Range("F21").Value = 1.039168003931558
$ws.Range("I21").Value = 1.030618562053732
$ws.Range("J21").Value = 1.02943714188261
$ws.Range("K21").Value = 1.035141868089567
$ws.Range("L21").Value = 1.027151779324125
$ws.Range("M21").Value = 1.042736391765126
$ws.Range("N21").Value = 1.013739526591588
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021936663537913
$ws.Range("D22").Value = 1.030923974026713
$ws.Range("E22").Value = 1.022828121205451
$ws.Range("F22").Value = 1.038380994760536
$ws.Range("I22").Value = 1.030472943644543
$ws.Range("J22").Value = 1.028873405143244
$ws.Range("K22").Value = 1.03465281785889
$ws.Range("L22").Value = 1.026588875259131
$ws.Range("M22").Value = 1.042080935017846
$ws.Range("N22").Value = 1.013547558476881
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022370232104972
$ws.Range("D23").Value = 1.031253249642531
$ws.Range("E23").Value = 1.023197296838555
$ws.Range("F23").Value = 1.038798177019147
$ws.Range("I23").Value = 1.030550286790301
$ws.Range("J23").Value = 1.02917230529928
$ws.Range("K23").Value = 1.034912157278936
$ws.Range("L23").Value = 1.026887309603895
$ws.Range("M23").Value = 1.042428437514916
$ws.Range("N23").Value = 1.013649357551258
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024077055205866
$ws.Range("D24").Value = 1.032549129235453
$ws.Range("E24").Value = 1.024651446525807
$ws.Range("F24").Value = 1.040441169321743
$ws.Range("I24").Value = 1.030851496041473
$ws.Range("J24").Value = 1.030347883489485
$ws.Range("K24").Value = 1.035931277385825
$ws.Range("L24").Value = 1.028061594295834
$ws.Range("M24").Value = 1.043795805384076
$ws.Range("N24").Value = 1.014049396637835
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026057583428999
$ws.Range("D25").Value = 1.034052003736867
$ws.Range("E25").Value = 1.026340454114777
$ws.Range("F25").Value = 1.042349003016116
$ws.Range("I25").Value = 1.031194231729649
$ws.Range("J25").Value = 1.031709686750222
$ws.Range("K25").Value = 1.037110036628013
$ws.Range("L25").Value = 1.029423011867182
$ws.Range("M25").Value = 1.045381091887723
$ws.Range("N25").Value = 1.014512102007575
